$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.94"
$ws.Range("E2").Value = "'5.06%"
$ws.Range("D3").Value = "'44.53"
$ws.Range("E3").Value = "'7.24%"
$ws.Range("D4").Value = "'5.099"
$ws.Range("E4").Value = "'1.29%"
$ws.Range("D5").Value = "'0.08014"
$ws.Range("E5").Value = "'6.12%"
$ws.Range("D6").Value = "'4.484"
$ws.Range("E6").Value = "'2.40%"
$ws.Range("D7").Value = "'1.650"
$ws.Range("E7").Value = "'3.83%"
$ws.Range("D8").Value = "'1.079"
$ws.Range("E8").Value = "'16.30%"
$ws.Range("D9").Value = "'0.1288"
$ws.Range("E9").Value = "'7.00%"
$ws.Range("D10").Value = "'0.1894"
$ws.Range("E10").Value = "'4.03%"
$ws.Range("D11").Value = "'0.09232"
$ws.Range("E11").Value = "'3.51%"
$ws.Range("D12").Value = "'0.04199"
$ws.Range("E12").Value = "'7.01%"
$ws.Range("D13").Value = "'0.1037"
$ws.Range("E13").Value = "'-1.67%"
$ws.Range("D14").Value = "'0.001306"
$ws.Range("E14").Value = "'1.13%"
$ws.Range("D15").Value = "'0.005848"
$ws.Range("E15").Value = "'-0.45%"
$ws.Range("D17").Value = "'3.372"
$ws.Range("E17").Value = "'1.16%"
$ws.Range("D18").Value = "'2.402"
$ws.Range("E18").Value = "'-0.97%"
$ws.Range("D19").Value = "'0.3361"
$ws.Range("E19").Value = "'1.30%"
$ws.Range("D20").Value = "'7.995"
$ws.Range("E20").Value = "'0.51%"
$ws.Range("D21").Value = "'0.1378"
$ws.Range("E21").Value = "'-2.93%"
$ws.Range("D22").Value = "'0.3129"
$ws.Range("E22").Value = "'4.32%"
$ws.Range("D23").Value = "'0.04192"
$ws.Range("E23").Value = "'3.31%"
$ws.Range("E24").Value = "'0.26%"
$ws.Range("D25").Value = "'0.004602"
$ws.Range("E25").Value = "'15.23%"
$ws.Range("D26").Value = "'0.0001336"
$ws.Range("E26").Value = "'8.41%"
$ws.Range("E38").Value = "'10.28%"
$ws.Range("D39").Value = "'0.05419"
$ws.Range("E39").Value = "'3.83%"
$ws.Range("D40").Value = "'0.005607"
$ws.Range("E40").Value = "'-12.40%"
$ws.Range("D41").Value = "'0.007717"
$ws.Range("E41").Value = "'-0.85%"
$ws.Range("D42").Value = "'0.1412"
$ws.Range("E42").Value = "'6.23%"
$ws.Range("D43").Value = "'0.007299"
$ws.Range("E43").Value = "'-3.80%"
$ws.Range("D44").Value = "'0.008387"
$ws.Range("E44").Value = "'6.87%"
$ws.Range("D45").Value = "'0.3122"
$ws.Range("E45").Value = "'-3.01%"
$ws.Range("D46").Value = "'0.00006703"
$ws.Range("E46").Value = "'-1.28%"
$ws.Range("E47").Value = "'-1.20%"
$ws.Range("D48").Value = "'0.05303"
$ws.Range("E48").Value = "'15.14%"
$ws.Range("D49").Value = "'0.003958"
$ws.Range("E49").Value = "'-5.91%"
$ws.Range("D50").Value = "'0.00002078"
$ws.Range("E50").Value = "'-1.20%"
$ws.Range("D51").Value = "'0.0001979"
$ws.Range("E51").Value = "'-1.20%"
